# ---------------------------------------------------------------------------
# attendance.xlsx — "Add files via upload" edit
#
# Summary of the target state (per the OOXML diff):
#   1. Insert a new empty sheet "ListOfManagers" between "ListOfTrainers" and
#      "CourseDetails".
#   2. Insert a new sheet "1212_1212" at the very end, containing an
#      attendance sheet export.
#   3. Rewrite the data on ListOfTrainees, ListOfTrainers, CourseDetails,
#      MappingCourseTrainees and MappingCourseTrainers.
#   4. Fix up selections / active sheet so the final active tab is
#      "MappingCourseTrainers".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet structure: insert the two new sheets at the right spots.
# ---------------------------------------------------------------------------

$wsTrainers = $wb.Worksheets.Item("ListOfTrainers")
$wsManagers = $wb.Worksheets.Add($null, $wsTrainers)
$wsManagers.Name = "ListOfManagers"
$wsManagers.Outline.SummaryRow = 1
$wsManagers.Outline.SummaryColumn = 1

$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAttendance = $wb.Worksheets.Add($null, $wsLast)
$wsAttendance.Name = "1212_1212"
$wsAttendance.Outline.SummaryRow = 1
$wsAttendance.Outline.SummaryColumn = 1

# ---------------------------------------------------------------------------
# 2. ListOfTrainees
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("ListOfTrainees")
$ws1.Cells.ClearContents()

$ws1.Range("A1").Value = 1
$ws1.Range("B1").Value = "John Wick"
$ws1.Range("C1").Value = "C343"
$ws1.Range("D1").Value = "degree"

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Caine"
$ws1.Range("C3").Value = "C343"
$ws1.Range("D3").Value = "work"

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "killa"
$ws1.Range("C4").Value = "C343"
$ws1.Range("D4").Value = "dksdf"

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "ksdf"
$ws1.Range("C5").Value = "C333"
$ws1.Range("D5").Value = "degree"

[void]$ws1.Range("B1").Select()

# ---------------------------------------------------------------------------
# 3. ListOfTrainers
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Item("ListOfTrainers")
$ws2.Cells.ClearContents()

$ws2.Range("A1").Value = 1
$ws2.Range("B1").Value = "A test"
$ws2.Range("C1").Value = "sf"
$ws2.Range("D1").Value = 33
$ws2.Range("E1").Value = "C343"

[void]$ws2.Range("B1").Select()

# ---------------------------------------------------------------------------
# 4. ListOfManagers (brand-new, stays empty)
# ---------------------------------------------------------------------------

[void]$wsManagers.Range("I34").Select()

# ---------------------------------------------------------------------------
# 5. CourseDetails
# ---------------------------------------------------------------------------

$ws4 = $wb.Worksheets.Item("CourseDetails")
$ws4.Cells.ClearContents()

$ws4.Range("A1").Value = "C343"
$ws4.Range("B1").Value = "Python and C++"

[void]$ws4.Range("B1").Select()

# ---------------------------------------------------------------------------
# 6. MappingCourseTrainees
# ---------------------------------------------------------------------------

$ws5 = $wb.Worksheets.Item("MappingCourseTrainees")
$ws5.Cells.ClearContents()

$ws5.Range("A1").Value = "C343"
$ws5.Range("B1").Value = "1, 2, 3"

[void]$ws5.Range("B1").Select()

# ---------------------------------------------------------------------------
# 7. "1212_1212" — attendance export sheet
# ---------------------------------------------------------------------------

$ws7 = $wsAttendance

$ws7.Range("A1").Value = "12:30"
$ws7.Range("B1").Value = "16:00"

$ws7.Range("A2").Value = "Trainer"
$ws7.Range("B2").Value = "A test"

$ws7.Range("A3").Value = "Attendance of trainees"

$ws7.Range("A4").Value = 1
$ws7.Range("B4").Value = "John Wick"
$ws7.Range("C4").Value = "A"

$ws7.Range("A6").Value = 2
$ws7.Range("B6").Value = "Caine"
$ws7.Range("C6").Value = "P"

$ws7.Range("A7").Value = 3
$ws7.Range("B7").Value = "killa"
$ws7.Range("C7").Value = "A"

[void]$ws7.Range("A1").Select()

# ---------------------------------------------------------------------------
# 8. MappingCourseTrainers — written/activated LAST so it ends up the
#    workbook's active tab/sheet on save (matches activeTab="5" in the diff).
# ---------------------------------------------------------------------------

$ws6 = $wb.Worksheets.Item("MappingCourseTrainers")
$ws6.Cells.ClearContents()

$ws6.Range("A1").Value = "C343"

# B1 must be the literal text "1" (not the number 1) per the target file.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "1"

# This sheet ends up being the active tab/sheet in the saved workbook.
[void]$ws6.Activate()
[void]$ws6.Range("A1").Select()
